# Update the "Förändrad" (changed) date column for all existing data rows
# (rows 2-171) from 45205 (2023-10-06) to 45206 (2023-10-07).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C171").Value = 45206

# Row 171 gains an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(171).RowHeight = 15

# Append the new record as row 172.
$ws.Range("A172").Value = "A 48047-2023"
$ws.Range("B172").Value = 45204
$ws.Range("C172").Value = 45206
$ws.Range("B172:C172").NumberFormat = "YYYY-MM-DD"
$ws.Range("D172").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E172").Value = "VÄNNÄS"
$ws.Range("G172").Value = 16.4
$ws.Range("H172").Value = 0
$ws.Range("I172").Value = 0
$ws.Range("J172").Value = 0
$ws.Range("K172").Value = 0
$ws.Range("L172").Value = 0
$ws.Range("M172").Value = 0
$ws.Range("N172").Value = 0
$ws.Range("O172").Value = 0
$ws.Range("P172").Value = 0
$ws.Range("Q172").Value = 0
$ws.Range("R172").WrapText = $true
